$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B ("Year of Treatment"); remaining columns C..K shift left to B..J
$ws.Range("B1").EntireColumn.Delete()

# Append ".global" suffix to each header cell (now in B1:J1, previously C1:K1)
for ($col = 2; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value()
    $cell.Value = "$current.global"
}
